# Apply updated Kraken market-price / profit figures across all leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 345.35715
$ws.Range("J9").Value = 44.5
$ws.Range("L9").Value = 44.5
$ws.Range("N9").Value = -382.5

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0

$ws.Range("H33").Value = 116.42857
$ws.Range("I33").Value = 125.833336
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 125.833336
$ws.Range("L33").Value = 60
$ws.Range("M33").Value = 103.166664
$ws.Range("N33").Value = -518

$ws.Range("H88").Value = 3269.8
$ws.Range("I88").Value = 3337.25
$ws.Range("K88").Value = 3337.25
$ws.Range("M88").Value = -2931.25

$ws.Range("H91").Value = 3269.8
$ws.Range("I91").Value = 3337.25
$ws.Range("K91").Value = 3337.25
$ws.Range("M91").Value = -1933.25

$ws.Range("H103").Value = 665

$ws.Range("H137").Value = 1828.4445
$ws.Range("I137").Value = 1912.4286
$ws.Range("J137").Value = 1534.5
$ws.Range("K137").Value = 5737.2858
$ws.Range("L137").Value = 4603.5
$ws.Range("M137").Value = -3187.2858
$ws.Range("N137").Value = -9703.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 21999.666
$ws.Range("J92").Value = 21999.666
$ws.Range("L92").Value = 21999.666
$ws.Range("N92").Value = -26991.666

$ws.Range("H132").Value = 3330.3333
$ws.Range("I132").Value = 3396.4
$ws.Range("K132").Value = 10189.2
$ws.Range("M132").Value = -7659.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 37000
$ws.Range("J35").Value = 37000
$ws.Range("L35").Value = 37000
$ws.Range("N35").Value = -37620

$ws.Range("H105").Value = 4373.375
$ws.Range("I105").Value = 4197.4
$ws.Range("K105").Value = 4197.4
$ws.Range("M105").Value = -2450.4

$ws.Range("H134").Value = 9777.091
$ws.Range("J134").Value = 13508.5
$ws.Range("L134").Value = 40525.5
$ws.Range("N134").Value = -45595.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 3088.75
$ws.Range("J11").Value = 3088.75
$ws.Range("L11").Value = 3088.75
$ws.Range("N11").Value = -3368.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2763.6667
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H135").Value = 2763.6667
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5580
$ws.Range("I43").Value = 5580
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 5580
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -5429

$ws.Range("H46").Value = 12783.333

$ws.Range("H57").Value = 5000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H80").Value = 6309.625
$ws.Range("I80").Value = 6354
$ws.Range("J80").Value = 5999
$ws.Range("K80").Value = 6354
$ws.Range("L80").Value = 5999
$ws.Range("M80").Value = -5356
$ws.Range("N80").Value = -7995

$ws.Range("H83").Value = 6309.625
$ws.Range("I83").Value = 6354
$ws.Range("J83").Value = 5999
$ws.Range("K83").Value = 31770
$ws.Range("L83").Value = 29995
$ws.Range("M83").Value = -26778
$ws.Range("N83").Value = -39979

$ws.Range("H122").Value = 6799.5713
$ws.Range("I122").Value = 5516.3335
$ws.Range("K122").Value = 16549.0005
$ws.Range("M122").Value = -14099.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2947
$ws.Range("I7").Value = 3009.8462
$ws.Range("J7").Value = 2742.75
$ws.Range("K7").Value = 3009.8462
$ws.Range("L7").Value = 2742.75
$ws.Range("M7").Value = -2897.8462
$ws.Range("N7").Value = -2966.75

$ws.Range("H40").Value = 4734.9
$ws.Range("I40").Value = 4419.143
$ws.Range("J40").Value = 5471.6665
$ws.Range("K40").Value = 4419.143
$ws.Range("L40").Value = 5471.6665
$ws.Range("M40").Value = -4283.143
$ws.Range("N40").Value = -5743.6665

$ws.Range("H61").Value = 7500
$ws.Range("I61").Value = 7500
$ws.Range("K61").Value = 7500
$ws.Range("M61").Value = -7298

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("N69").Value = 0

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("N72").Value = 0

$ws.Range("H113").Value = 7500
$ws.Range("I113").Value = 7500
$ws.Range("K113").Value = 7500
$ws.Range("M113").Value = -5330

$ws.Range("H122").Value = 4667.1113
$ws.Range("I122").Value = 4375.5
$ws.Range("K122").Value = 13126.5
$ws.Range("M122").Value = -10676.5

$ws.Range("H126").Value = 2947
$ws.Range("I126").Value = 3009.8462
$ws.Range("J126").Value = 2742.75
$ws.Range("K126").Value = 9029.5386
$ws.Range("L126").Value = 8228.25
$ws.Range("M126").Value = -6559.5386
$ws.Range("N126").Value = -13168.25

$ws.Range("H132").Value = 3335.25
$ws.Range("I132").Value = 3335.25
$ws.Range("K132").Value = 10005.75
$ws.Range("M132").Value = -7475.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2523.2727
$ws.Range("I2").Value = 94.75
$ws.Range("J2").Value = 8999.333000000001
$ws.Range("K2").Value = 94.75
$ws.Range("L2").Value = 8999.333000000001
$ws.Range("M2").Value = 17.25
$ws.Range("N2").Value = -9223.333000000001

$ws.Range("H40").Value = 19990
$ws.Range("I40").Value = 19990
$ws.Range("K40").Value = 19990
$ws.Range("M40").Value = -19841

$ws.Range("H81").Value = 2266.6667
$ws.Range("I81").Value = 2400
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 4800
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -3739
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 2266.6667
$ws.Range("I84").Value = 2400
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 24000
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -18696
$ws.Range("N84").Value = -30608

$ws.Range("H132").Value = 3712.2666
$ws.Range("I132").Value = 1517
$ws.Range("K132").Value = 4551
$ws.Range("M132").Value = -2021
